# Update the build version string across the workbook to reflect the new
# release build timestamp.
#
# Old: mines - January 30 (built on January 30 2026 16.19.47 EST)
# New: mines - January 30 (built on February 02 2026 12.49.33 EST)

$wb = $excel.ActiveWorkbook

$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet  = $wb.Worksheets.Item("Boundaries and methane sources")

# "About" sheet: standalone version line and the citation sentence that
# embeds the version string.
$aboutSheet.Range("A2").Value = "Version: " + $newVersion

$citationNew = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Fangezhuang Coal Mine, China, M1121, version ' + "'" + $newVersion + "'" + '. (See the CC license for attribution requirements if sharing or adapting the data set.)'
$aboutSheet.Range("A6").Value = $citationNew

# "Boundaries and methane sources" sheet: build_version column (S) for each
# data row (rows 2 through 8).
for ($row = 2; $row -le 8; $row++) {
    $dataSheet.Cells.Item($row, 19).Value = $newVersion
}
